# Weekly CompStat crime report (Brooklyn South) refresh.
# New crime data collected for the week of 9/4/2023 through 9/10/2023
# (previously the report covered 8/28/2023 through 9/3/2023), with the
# report volume/number bumped from "35" to "36".
#
# The weekly "Crime Complaints" table (rows 14-30) and the historical
# "2 Year" column (row 35, J35) stay the same shape; only the reported
# figures are refreshed with the newly collected counts/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Murder (row 14) ---
# C14 was a blank-week placeholder (text "0"); now a real reported count,
# so give it the same numeric "#,##0" look as the rest of the row.
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = -50
$ws.Range("I14").Value = 42
$ws.Range("J14").Value = 43
$ws.Range("K14").Value = -2.325581395348
$ws.Range("L14").Value = 27.272727272727
$ws.Range("M14").Value = -34.375
$ws.Range("N14").Value = -75.438596491228

# --- Rape (row 15) ---
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 200
$ws.Range("G15").Value = 14
$ws.Range("H15").Value = 28.571428571428
$ws.Range("I15").Value = 151
$ws.Range("J15").Value = 147
$ws.Range("K15").Value = 2.721088435374
$ws.Range("L15").Value = 0.666666666666
$ws.Range("M15").Value = 17.96875
$ws.Range("N15").Value = -62.899262899262

# --- Robbery (row 16) ---
$ws.Range("C16").Value = 29
$ws.Range("D16").Value = 41
$ws.Range("E16").Value = -29.268292682926
$ws.Range("F16").Value = 126
$ws.Range("G16").Value = 149
$ws.Range("H16").Value = -15.436241610738
$ws.Range("I16").Value = 1210
$ws.Range("J16").Value = 1356
$ws.Range("K16").Value = -10.766961651917
$ws.Range("L16").Value = 32.096069868995
$ws.Range("M16").Value = -37.948717948717
$ws.Range("N16").Value = -87.313902285594

# --- Fel. Assault (row 17) ---
$ws.Range("C17").Value = 80
$ws.Range("D17").Value = 83
$ws.Range("E17").Value = -3.614457831325
$ws.Range("F17").Value = 287
$ws.Range("G17").Value = 315
$ws.Range("H17").Value = -8.888888888888
$ws.Range("I17").Value = 2479
$ws.Range("J17").Value = 2441
$ws.Range("K17").Value = 1.556739041376
$ws.Range("L17").Value = 18.953934740882
$ws.Range("M17").Value = 43.129330254041
$ws.Range("N17").Value = -48.854961832061

# --- Burglary (row 18) ---
$ws.Range("C18").Value = 30
$ws.Range("E18").Value = -30.232558139534
$ws.Range("F18").Value = 158
$ws.Range("G18").Value = 176
$ws.Range("H18").Value = -10.227272727272
$ws.Range("I18").Value = 1218
$ws.Range("J18").Value = 1456
$ws.Range("K18").Value = -16.346153846153
$ws.Range("L18").Value = 6.654991243432
$ws.Range("M18").Value = -48.607594936708
$ws.Range("N18").Value = -90.42076287849

# --- Gr. Larceny (row 19) ---
$ws.Range("C19").Value = 103
$ws.Range("D19").Value = 130
$ws.Range("E19").Value = -20.76923076923
$ws.Range("F19").Value = 481
$ws.Range("G19").Value = 591
$ws.Range("H19").Value = -18.612521150592
$ws.Range("I19").Value = 4553
$ws.Range("J19").Value = 4993
$ws.Range("K19").Value = -8.812337272181
$ws.Range("L19").Value = 36.891160553217
$ws.Range("M19").Value = 21.058229194363
$ws.Range("N19").Value = -25.213534822601

# --- G.L.A. (row 20) ---
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 40
$ws.Range("E20").Value = -7.5
$ws.Range("F20").Value = 156
$ws.Range("G20").Value = 137
$ws.Range("H20").Value = 13.868613138686
$ws.Range("I20").Value = 1288
$ws.Range("J20").Value = 1264
$ws.Range("K20").Value = 1.898734177215
$ws.Range("L20").Value = 49.593495934959
$ws.Range("M20").Value = -6.053975200583
$ws.Range("N20").Value = -92.029702970297

# --- TOTAL (row 21) ---
$ws.Range("C21").Value = 286
$ws.Range("D21").Value = 340
$ws.Range("E21").Value = -15.882352941176
$ws.Range("F21").Value = 1229
$ws.Range("G21").Value = 1388
$ws.Range("H21").Value = -11.455331412103
$ws.Range("I21").Value = 10941
$ws.Range("J21").Value = 11700
$ws.Range("K21").Value = -6.487179487179
$ws.Range("L21").Value = 28.536184210526
$ws.Range("M21").Value = -3.823839662447
$ws.Range("N21").Value = -78.085566638625

# --- Transit (row 22) ---
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 128
$ws.Range("J22").Value = 132
$ws.Range("K22").Value = -3.030303030303
$ws.Range("L22").Value = 29.292929292929
$ws.Range("M22").Value = -37.560975609756

# --- Housing (row 23) ---
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -14.285714285714
$ws.Range("F23").Value = 45
$ws.Range("G23").Value = 49
$ws.Range("H23").Value = -8.163265306122
$ws.Range("I23").Value = 373
$ws.Range("J23").Value = 387
$ws.Range("K23").Value = -3.617571059431
$ws.Range("L23").Value = 12.349397590361
$ws.Range("M23").Value = 54.132231404958

# --- Petit Larceny (row 24) ---
$ws.Range("C24").Value = 288
$ws.Range("D24").Value = 309
$ws.Range("E24").Value = -6.796116504854
$ws.Range("F24").Value = 1226
$ws.Range("G24").Value = 1416
$ws.Range("H24").Value = -13.418079096045
$ws.Range("I24").Value = 11103
$ws.Range("J24").Value = 11130
$ws.Range("K24").Value = -0.242587601078
$ws.Range("L24").Value = 40.918898337352
$ws.Range("M24").Value = 30.072633552015

# --- Misd. Assault (row 25) ---
$ws.Range("C25").Value = 123
$ws.Range("D25").Value = 126
$ws.Range("E25").Value = -2.380952380952
$ws.Range("F25").Value = 457
$ws.Range("G25").Value = 447
$ws.Range("H25").Value = 2.237136465324
$ws.Range("I25").Value = 4099
$ws.Range("J25").Value = 3919
$ws.Range("K25").Value = 4.593008420515
$ws.Range("L25").Value = 24.438372799028
$ws.Range("M25").Value = -13.413603717786

# --- UCR Rape* (row 26) ---
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 60
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 4.166666666666
$ws.Range("I26").Value = 218
$ws.Range("J26").Value = 232
$ws.Range("K26").Value = -6.03448275862
$ws.Range("L26").Value = -7.627118644067

# --- Other Sex Crimes (row 27) ---
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = -11.111111111111
$ws.Range("F27").Value = 43
$ws.Range("G27").Value = 46
$ws.Range("H27").Value = -6.521739130434
$ws.Range("I27").Value = 460
$ws.Range("J27").Value = 481
$ws.Range("K27").Value = -4.365904365904
$ws.Range("L27").Value = 9.263657957244

# --- Shooting Vic. (row 28) ---
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = -42.857142857142
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 23
$ws.Range("H28").Value = -56.521739130434
$ws.Range("I28").Value = 106
$ws.Range("J28").Value = 159
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -23.741007194244
$ws.Range("M28").Value = -46.192893401015
$ws.Range("N28").Value = -81.138790035587

# --- Shooting Inc. (row 29) ---
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -20
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -37.5
$ws.Range("I29").Value = 93
$ws.Range("J29").Value = 121
$ws.Range("K29").Value = -23.140495867768
$ws.Range("L29").Value = -24.390243902439
$ws.Range("M29").Value = -42.23602484472
$ws.Range("N29").Value = -81.097560975609

# --- Hate Crimes (row 30) ---
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 52
$ws.Range("K30").Value = -42.222222222222
$ws.Range("L30").Value = 10.63829787234
